# Adapt column header formatting to respective input file names (#7)
#
# - Header cells that previously carried a "_old" suffix now carry a
#   "_FV2404" suffix, and header cells that previously carried a "_new"
#   suffix now carry a "_FV2410" suffix (the "diff" header is untouched).
# - The header row is frozen (pane split after row 1).
# - The whole sheet range is wrapped in a native Excel table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------
# Columns A..J used to end in "_old"  -> now end in "_FV2404"
# Column  K is "diff"                 -> stays "diff"
# Columns L..U used to end in "_new"  -> now end in "_FV2410"
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2404"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2410"
}

# --- 2. Freeze the header row --------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the data range in a table ------------------------------------
$tableRange = $ws.Range("A1:U66")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
